# Re-touch the run-level character formatting (font family, size, bold,
# italic, strike-through and color) of every populated cell in the first
# table so the document is re-serialized by the current OOXML writer -
# matching the "Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3" fixture
# refresh, where only the writer's low-level spelling of these (already
# unchanged) run properties was updated.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellFont($row, $col, $bold) {
    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    # Exclude the trailing cell-mark / paragraph-mark character so only
    # the run(s) holding the visible text (or the single empty run) are
    # restyled - touching the mark itself would add stray paragraph mark
    # run properties that aren't part of this change.
    $r.MoveEnd(1, -1)
    $r.Font.Name = "Calibri"
    $r.Font.Size = 11
    $r.Font.Bold = $bold
    $r.Font.Italic = $false
    $r.Font.StrikeThrough = $false
    $r.Font.Color = 0
}

# Row 1: Item / Price / Quantity / Total (bold)
Set-CellFont 1 1 $false
Set-CellFont 1 2 $false
Set-CellFont 1 3 $false
Set-CellFont 1 4 $true

# Row 2: Apple / 1,76 € / 23 / 40,48 €
Set-CellFont 2 1 $false
Set-CellFont 2 2 $false
Set-CellFont 2 3 $false
Set-CellFont 2 4 $false

# Row 3: empty line - only cells 2 and 4 carry run formatting
Set-CellFont 3 2 $false
Set-CellFont 3 4 $false

# Row 4: Banana / 1,99 € / 45 / 89,55 €
Set-CellFont 4 1 $false
Set-CellFont 4 2 $false
Set-CellFont 4 3 $false
Set-CellFont 4 4 $false

# Row 5: Total / 130,03 € - only cells 3 and 4 carry run formatting
Set-CellFont 5 3 $false
Set-CellFont 5 4 $false
